$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "01"-suffixed header labels in row 1 one column to the right (H1..O1 -> I1..P1),
# clearing H1 entirely (copy the blank/default formatting from P1 so the cell serializes as empty).
$ws.Range("P1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").ClearContents()

# P1 inherits the column's normal style (like A1) now that it holds real text.
$ws.Range("A1").Copy()
$ws.Range("P1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "ADTYMT01"
$ws.Range("J1").Value = "GDRID101"
$ws.Range("K1").Value = "W08PPR01"
$ws.Range("L1").Value = "LMSACN01"
$ws.Range("M1").Value = "LMSAPN01"
$ws.Range("N1").Value = "W08LBL01"
$ws.Range("O1").Value = "W08DLY01"
$ws.Range("P1").Value = "ACTACT"

# Row 2: swap the period labels and populate the new "戶況" column.
$ws.Range("A2").Value = "前期年月份"
$ws.Range("H2").Value = "戶況"
$ws.Range("I2").Value = "當期年月份"

$ws.Range("B15").Select()
